# B6-PowerPoint.pptx edit — Tue, Apr 21, 2020  1:05:40 PM
#
# 1. The presentation's colour theme is switched from the "Integral" /
#    "Red Violet" palette over to the plain built-in "Office" palette
#    (this is what ends up serialised into ppt/theme/theme1.xml).
# 2. The three data tables (on the slides that used the bespoke
#    "Table_0" style) are switched over to PowerPoint's built-in
#    "No Style, Table Grid" table style.

$p = $ppt.ActivePresentation

function RGBVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# --- 1. Swap the colour scheme over to the default Office palette ----------
# Order matches the standard ppColorSchemeIndex layout:
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1..accent6, 11 hlink, 12 folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = RGBVal $officeColors[$i - 1]
}

# --- 2. Re-style the three tables with the built-in "No Style, Table Grid" -
$newTableStyle = "{B84E7C2A-A07A-4224-A9F3-DB3372044DBA}"

for ($slideIdx = 14; $slideIdx -le 16; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
